# fichier_usager_test.xlsx — "make nir import more secure"
#
# 1) The date-of-birth sample value in E2 is adjusted (year corrected).
# 2) A batch of cell formats that previously had no fill now pick up the
#    same light fill (fillId 4 == indexed color 11) that the rest of the
#    sheet's "bordered" formats already use, so the whole grid reads as
#    one consistent banded style instead of a mix of filled/unfilled
#    cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) fix the sample date of birth in E2 -------------------------------
$ws.Range("E2").Value = 29546

# --- 2) add the missing fill to the previously "no fill" formats --------
# Column E (date of birth) rows 2-8 (numFmtId 59 date format)
$ws.Range("E2:E8").Interior.ColorIndex = 4

# Column I row 2 ("Role" column header group)
$ws.Range("I2").Interior.ColorIndex = 4

# One-off cells that picked up their own tweaked format (rows 3-5)
$ws.Range("C3").Interior.ColorIndex = 4
$ws.Range("P3").Interior.ColorIndex = 4
$ws.Range("P4").Interior.ColorIndex = 4
$ws.Range("F5").Interior.ColorIndex = 4
$ws.Range("G5").Interior.ColorIndex = 4

# Data block rows 9-67 : column J separately from A-I / K-P
# (the engine only honours the first area of a multi-area Range, so each
# contiguous block is set individually rather than via a comma union)
$ws.Range("J9:J67").Interior.ColorIndex = 4
$ws.Range("A37:I67").Interior.ColorIndex = 4
$ws.Range("K9:P67").Interior.ColorIndex = 4

# Subtotal row 68
$ws.Range("A68:I68").Interior.ColorIndex = 4
$ws.Range("K68:P68").Interior.ColorIndex = 4
$ws.Range("J68").Interior.ColorIndex = 4

# Remaining data block rows 69-220
$ws.Range("A69:I220").Interior.ColorIndex = 4
$ws.Range("K69:P220").Interior.ColorIndex = 4
$ws.Range("J69:J220").Interior.ColorIndex = 4
